# Applies the "Need refuse to answer option?" column (N) addition plus the
# "Neutral" highlight formatting on several Field-Name cells in column A,
# as described by the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New column N: header text + width
# ---------------------------------------------------------------------
$ws.Range("N3").Value = "Need refuse to answer option?"

# Target stored column width is 27.5 characters; this runtime quantizes
# ColumnWidth to an internal pixel grid (Calibri 11 / MDW=7), so we feed it
# the input that rounds to the closest achievable stored width (27.428571...).
$ws.Columns.Item(14).ColumnWidth = 26.7142857142857

# ---------------------------------------------------------------------
# New cell N6 ("Yes", centered like the other Yes/No flags in that row)
# ---------------------------------------------------------------------
$ws.Range("N6").Value = "Yes"
$ws.Range("N6").HorizontalAlignment = -4108   # xlCenter

# Reflect the new active selection left behind by the edit
[void]$ws.Range("N6").Select()

# ---------------------------------------------------------------------
# Rows whose Field Name (column A) gets the built-in "Neutral" highlight
# style, along with the taller (15pt) row height that comes with it
# ---------------------------------------------------------------------
$rowsWithNeutralFieldName = @(7, 8, 12, 15, 16, 20, 21, 22, 23, 24, 25, 26, 27, 29)
foreach ($r in $rowsWithNeutralFieldName) {
    $ws.Range("A$r").Style = "Neutral"
    $ws.Rows.Item($r).RowHeight = 15
}

# ---------------------------------------------------------------------
# Existing (previously blank) column-N cells that now get "Yes"
# ---------------------------------------------------------------------
$rowsWithNewYesFlag = @(19, 20, 21, 22, 23, 24, 27, 28, 30, 31, 32, 33, 36, 37)
foreach ($r in $rowsWithNewYesFlag) {
    $ws.Range("N$r").Value = "Yes"
}
